$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price + 1h volume change) plus two
# coin re-rankings (rows 10/11 and rows 39/40 swapped places).

$ws.Range("D2").Value = "64.124.43"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.760.39"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'577.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'159.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.25%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.93%  "
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "3.251.01"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "'27.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "63.734.27"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "2.764.04"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "'359.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'65.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").Value = "'0.172"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "0.0₃0917"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'7.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "'1.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "'169.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'4.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'351.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'39.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'21.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "'22.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").Value = "'0.0593"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").Value = "'137.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "'0.635"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'11.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
